$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 28.95628266666667
$ws.Range("H2").Value = 86.868848
$ws.Range("I2").Value = 0.5491054194301004
$ws.Range("J2").Value = 0.5491054194301005
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.916802666666667
$ws.Range("N2").Value = 5.750408
$ws.Range("O2").Value = 0.3865473586068074
$ws.Range("P2").Value = 0.3865473586068074
$ws.Range("Q2").Value = 55.50347983222044
$ws.Range("R2").Value = 499.531318489984
$ws.Range("S2").Value = 0.2122552494773884
$ws.Range("T2").Value = 0.2122552494773885

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 28.95628266666667
$ws.Range("H3").Value = 86.868848
$ws.Range("I3").Value = 0.5491054194301004
$ws.Range("J3").Value = 0.5491054194301005
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.041975666666666
$ws.Range("N3").Value = 9.125926999999999
$ws.Range("O3").Value = 0.6134526413931926
$ws.Range("P3").Value = 0.6134526413931926
$ws.Range("Q3").Value = 88.08430726912177
$ws.Range("R3").Value = 792.758765422096
$ws.Range("S3").Value = 0.3368501699527119
$ws.Range("T3").Value = 0.3368501699527121

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.691493
$ws.Range("H4").Value = 38.074479
$ws.Range("I4").Value = 0.2406720388519202
$ws.Range("J4").Value = 0.2406720388519202
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.916802666666667
$ws.Range("N4").Value = 5.750408
$ws.Range("O4").Value = 0.3865473586068074
$ws.Range("P4").Value = 0.3865473586068074
$ws.Range("Q4").Value = 24.32708762638133
$ws.Range("R4").Value = 218.943788637432
$ws.Range("S4").Value = 0.09303114090872466
$ws.Range("T4").Value = 0.09303114090872468

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.691493
$ws.Range("H5").Value = 38.074479
$ws.Range("I5").Value = 0.2406720388519202
$ws.Range("J5").Value = 0.2406720388519202
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.041975666666666
$ws.Range("N5").Value = 9.125926999999999
$ws.Range("O5").Value = 0.6134526413931926
$ws.Range("P5").Value = 0.6134526413931926
$ws.Range("Q5").Value = 38.60721287967033
$ws.Range("R5").Value = 347.464915917033
$ws.Range("S5").Value = 0.1476408979431955
$ws.Range("T5").Value = 0.1476408979431955

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4888703333333334
$ws.Range("H6").Value = 1.466611
$ws.Range("I6").Value = 0.009270573592685367
$ws.Range("J6").Value = 0.009270573592685367
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.916802666666667
$ws.Range("N6").Value = 5.750408
$ws.Range("O6").Value = 0.3865473586068074
$ws.Range("P6").Value = 0.3865473586068074
$ws.Range("Q6").Value = 0.9370679585875555
$ws.Range("R6").Value = 8.433611627288
$ws.Range("S6").Value = 0.003583515735022549
$ws.Range("T6").Value = 0.00358351573502255

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4888703333333334
$ws.Range("H7").Value = 1.466611
$ws.Range("I7").Value = 0.009270573592685367
$ws.Range("J7").Value = 0.009270573592685367
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.041975666666666
$ws.Range("N7").Value = 9.125926999999999
$ws.Range("O7").Value = 0.6134526413931926
$ws.Range("P7").Value = 0.6134526413931926
$ws.Range("Q7").Value = 1.487131658155222
$ws.Range("R7").Value = 13.384184923397
$ws.Range("S7").Value = 0.005687057857662817
$ws.Range("T7").Value = 0.005687057857662817

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Ephb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.59691233333333
$ws.Range("H8").Value = 31.790737
$ws.Range("I8").Value = 0.200951968125294
$ws.Range("J8").Value = 0.200951968125294
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.916802666666667
$ws.Range("N8").Value = 5.750408
$ws.Range("O8").Value = 0.3865473586068074
$ws.Range("P8").Value = 0.3865473586068074
$ws.Range("Q8").Value = 20.31218981896622
$ws.Range("R8").Value = 182.809708370696
$ws.Range("S8").Value = 0.07767745248567176
$ws.Range("T8").Value = 0.07767745248567177

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Ephb1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.59691233333333
$ws.Range("H9").Value = 31.790737
$ws.Range("I9").Value = 0.200951968125294
$ws.Range("J9").Value = 0.200951968125294
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.041975666666666
$ws.Range("N9").Value = 9.125926999999999
$ws.Range("O9").Value = 0.6134526413931926
$ws.Range("P9").Value = 0.6134526413931926
$ws.Range("Q9").Value = 32.23554945979988
$ws.Range("R9").Value = 290.1199451381989
$ws.Range("S9").Value = 0.1232745156396223
$ws.Range("T9").Value = 0.1232745156396223
